$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fully clear the old bold header row (content + formatting); it's being replaced ---
$ws.Range("A2:I2").Clear()

# --- Fully clear cells that are not part of the new layout at all ---
$ws.Range("G3").Clear()
$ws.Range("H3").Clear()
$ws.Range("I3").Clear()
$ws.Range("H4").Clear()
$ws.Range("G5").Clear()
$ws.Range("H5").Clear()
$ws.Range("D6").Clear()
$ws.Range("E6").Clear()
$ws.Range("G6").Clear()

# --- Write the new cell values (re-pivoted flavor grouping data) ---
$ws.Range("A1").Value = "Vanilla"
$ws.Range("B1").Value = "Caramel"
$ws.Range("C1").Value = "Toffee"
$ws.Range("D1").Value = "Honey"
$ws.Range("E1").Value = "Vanilla"
$ws.Range("A2").Value = "Fruity"
$ws.Range("B2").Value = "Dark Fruit"
$ws.Range("C2").Value = "Light Fruit"
$ws.Range("D2").Value = "Citrus"
$ws.Range("E2").Value = "Sour"
$ws.Range("F2").Value = "Fruity"
$ws.Range("A3").Value = "Earthy"
$ws.Range("B3").Value = "Coffee"
$ws.Range("C3").Value = "Chocolate"
$ws.Range("D3").Value = "Bitter Chocolate"
$ws.Range("E3").Value = "Nutty"
$ws.Range("F3").Value = "Earthy"
$ws.Range("A4").Value = "Mineral"
$ws.Range("B4").Value = "Dry"
$ws.Range("C4").Value = "Crisp"
$ws.Range("D4").Value = "Refreshing"
$ws.Range("E4").Value = "Clean"
$ws.Range("F4").Value = "Carbonation"
$ws.Range("G4").Value = "Mineral"
$ws.Range("A5").Value = "Toasty"
$ws.Range("B5").Value = "Roasty"
$ws.Range("C5").Value = "Rustic"
$ws.Range("D5").Value = "Smokey"
$ws.Range("E5").Value = "Toast"
$ws.Range("F5").Value = "Toasty"
$ws.Range("A6").Value = "Spices"
$ws.Range("B6").Value = "Spices"
$ws.Range("C6").Value = "Pepper"
$ws.Range("A7").Value = "Florals"
$ws.Range("B7").Value = "Herbal"
$ws.Range("C7").Value = "Grassy"
$ws.Range("D7").Value = "Hops"
$ws.Range("E7").Value = "Bitterness"
$ws.Range("F7").Value = "Florals"
$ws.Range("A8").Value = "Grainy"
$ws.Range("B8").Value = "Rye"
$ws.Range("C8").Value = "Corny"
$ws.Range("D8").Value = "Wheat"
$ws.Range("E8").Value = "Grainy"
$ws.Range("A9").Value = "Creamy"
$ws.Range("B9").Value = "Malt"
$ws.Range("C9").Value = "Creamy"

# --- Apply bold to column A (labels) for all data rows ---
$ws.Range("A1:A9").Font.Bold = $true

# --- Apply bold to the blank placeholder cells (creates them styled-but-empty) ---
$ws.Range("H2").Font.Bold = $true
$ws.Range("I2").Font.Bold = $true
$ws.Range("B13").Font.Bold = $true
$ws.Range("D13").Font.Bold = $true
$ws.Range("F13").Font.Bold = $true
$ws.Range("I13").Font.Bold = $true
$ws.Range("J13").Font.Bold = $true
$ws.Range("H14").Font.Bold = $true

# --- Column A: pin default width (matches the 10.83-char width in the source) ---
$ws.Columns.Item(1).ColumnWidth = 10

# --- Update selection to match the saved view state ---
$ws.Range("F3").Select()
